$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Rename the dataset headers (shared strings)
$ws.Range("B1").Value = "NUM_MERGES_PER_100_COMMITS_WO_FT"
$ws.Range("C1").Value = "NUM_MERGES_PER_100_COMMITS_WITH_FT"

# Resize columns B and C to fit the new, longer header text
$ws.Columns.Item(2).ColumnWidth = 39.833333333333336
$ws.Columns.Item(3).ColumnWidth = 41.5

# Select B1:C1 with active cell B1 (matches saved selection state)
$ws.Range("B1:C1").Select()
